# Add a new "About" sheet as the first tab in the workbook, describing
# where the framework's display Name is drawn from.
$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$aboutSheet = $wb.Worksheets.Add($firstSheet)
$aboutSheet.Name = "About"

# Header row
$aboutSheet.Range("A1").Value = "Name"
$aboutSheet.Range("B1").Value = "Description"

# Data row
$aboutSheet.Range("A2").Value = "Service"
$aboutSheet.Range("B2").Value = "Service Modalities cascade"

# Formatting: bold header, top-aligned data row
$aboutSheet.Range("A1:B1").Font.Bold = $true
$aboutSheet.Range("A2:B2").VerticalAlignment = -4160  # xlTop

# Leave selection on A3, as in the authored workbook
$aboutSheet.Range("A3").Select()
